# Update cryptos list values (Price and Volume(1h) columns) per scrape refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'51.792.06"
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").Value = "'2.779.91"
$ws.Range("E3").Value = '  -1.52%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'356.95"
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").Value = "'109.62"
$ws.Range("E6").Value = '  -1.84%  '
$ws.Range("E7").Value = '  -1.79%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("E9").Value = '  -1.91%  '
$ws.Range("D10").Value = "'39.84"
$ws.Range("E10").Value = '  -2.36%  '
$ws.Range("E11").Value = '  +2.42%  '
$ws.Range("E12").Value = '  -1.22%  '
$ws.Range("E13").Value = '  -2.24%  '
$ws.Range("D14").Value = "'7.61"
$ws.Range("E14").Value = '  -2.56%  '
$ws.Range("D15").Value = "'3.220.95"
$ws.Range("E15").Value = '  -1.24%  '
$ws.Range("D16").Value = "'2.782.80"
$ws.Range("E16").Value = '  -2.31%  '
$ws.Range("E17").Value = '  +1.68%  '
$ws.Range("D18").Value = "'51.710.63"
$ws.Range("E18").Value = '  -0.18%  '
$ws.Range("E19").Value = '  -1.95%  '
$ws.Range("E20").Value = '  -2.35%  '
$ws.Range("D21").Value = "'13.16"
$ws.Range("E21").Value = '  -1.84%  '
$ws.Range("D22").Value = "'0.0₃0970"
$ws.Range("E22").Value = '  -2.33%  '
$ws.Range("D23").Value = "'70.21"
$ws.Range("E23").Value = '  +0.19%  '
$ws.Range("D24").Value = "'269.33"
$ws.Range("D25").Value = "'2.74"
$ws.Range("E25").Value = '  -2.26%  '
$ws.Range("D26").Value = "'26.36"
$ws.Range("E26").Value = '  -2.49%  '
$ws.Range("D27").Value = "'1.00"
$ws.Range("E27").Value = '  -0.14%  '
$ws.Range("E28").Value = '  +16.72%  '
$ws.Range("D29").Value = "'10.23"
$ws.Range("E29").Value = '  -0.80%  '
$ws.Range("D30").Value = "'2.22"
$ws.Range("E30").Value = '  -1.56%  '
$ws.Range("D31").Value = "'6.26"
$ws.Range("E31").Value = '  +6.16%  '
$ws.Range("D32").Value = "'52.00"
$ws.Range("E32").Value = '  -1.02%  '
$ws.Range("D33").Value = "'34.69"
$ws.Range("E33").Value = '  -0.47%  '
$ws.Range("D34").Value = "'0.0450"
$ws.Range("E34").Value = '  -10.35%  '
$ws.Range("D35").Value = "'0.0840"
$ws.Range("E35").Value = '  -0.62%  '
$ws.Range("D36").Value = "'5.13"
$ws.Range("E36").Value = '  -6.10%  '
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("E38").Value = '  +1.85%  '
$ws.Range("E39").Value = '  -5.47%  '
$ws.Range("E40").Value = '  -3.81%  '
$ws.Range("D41").Value = "'2.55"
$ws.Range("E41").Value = '  +0.54%  '
$ws.Range("E42").Value = '  -2.31%  '
$ws.Range("E43").Value = '  -1.71%  '
$ws.Range("D44").Value = "'119.94"
$ws.Range("E44").Value = '  -4.69%  '
$ws.Range("D45").Value = "'21.83"
$ws.Range("E45").Value = '  -6.11%  '
$ws.Range("D46").Value = "'2.082.97"
$ws.Range("E46").Value = '  -0.75%  '
$ws.Range("D47").Value = "'3.26"
$ws.Range("E47").Value = '  -2.47%  '
$ws.Range("D48").Value = "'2.27"
$ws.Range("E48").Value = '  +0.58%  '
$ws.Range("D49").Value = "'5.74"
$ws.Range("E49").Value = '  -4.57%  '
$ws.Range("D50").Value = "'0.935"
$ws.Range("E50").Value = '  -5.56%  '
$ws.Range("E51").Value = '  +1.16%  '
